$wb = $excel.ActiveWorkbook

# Sheet "展览" - update "想去人数" (F) values for rows 3-5
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 3109
$ws1.Range("F4").Value = 146
$ws1.Range("F5").Value = 113

# Sheet "全部类型" - update "想去人数" (F) values for rows 7, 8, 10
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 3109
$ws4.Range("F8").Value = 146
$ws4.Range("F10").Value = 113
